$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23.. down by one
$ws.Rows.Item(23).Insert()

# Populate the new row 23 with the new weekly record
$ws.Cells.Item(23, 1).Value = 8
$ws.Cells.Item(23, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(23, 3).Value = "Coquimbo"
$ws.Cells.Item(23, 4).Value = 44972
$ws.Cells.Item(23, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23, 5).Value = 4
$ws.Cells.Item(23, 6).Value = 100114007
$ws.Cells.Item(23, 7).Value = "Jengibre"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 300
$ws.Cells.Item(23, 11).Value = 21000
$ws.Cells.Item(23, 12).Value = 22000
$ws.Cells.Item(23, 13).Value = 21500
$ws.Cells.Item(23, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(23, 15).Value = "Perú"
$ws.Cells.Item(23, 16).Value = 1654
$ws.Cells.Item(23, 17).Value = 13
$ws.Cells.Item(23, 18).Value = "Hortaliza"
